# Daily COVID stats refresh for the "Pais" worksheet (update to 18:17 snapshot).
# Mirrors the source diff: most rows just get refreshed B:H counters; three
# countries (Republica Dominicana, then Republica de Macedonia + Kenia) moved up
# in the case-count ranking, so the rows below them shift down with their
# untouched (pre-refresh) figures while the promoted countries get new numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 29 de Junio de 2020 a las 18:17'

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2650245
$ws.Range("C4").Value = 13168
$ws.Range("D4").Value = 1094372
$ws.Range("E4").Value = 1427360
$ws.Range("G4").Value = 76
$ws.Range("H4").Value = 128513

# Row 7: India
$ws.Range("B7").Value = 562457
$ws.Range("C7").Value = 13260
$ws.Range("D7").Value = 329728
$ws.Range("E7").Value = 215932
$ws.Range("G7").Value = 310
$ws.Range("H7").Value = 16797

# Row 8: Reino Unido
$ws.Range("B8").Value = 311965
$ws.Range("C8").Value = 814
$ws.Range("G8").Value = 25
$ws.Range("H8").Value = 43575

# Row 12: Italia
$ws.Range("B12").Value = 240436
$ws.Range("C12").Value = 126
$ws.Range("D12").Value = 189196
$ws.Range("E12").Value = 16496
$ws.Range("G12").Value = 6
$ws.Range("H12").Value = 34744

# Row 17: Alemania
$ws.Range("B17").Value = 195104
$ws.Range("C17").Value = 240
$ws.Range("E17").Value = 7974
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 9030

# Row 22: Canada
$ws.Range("B22").Value = 103818
$ws.Range("C22").Value = 568
$ws.Range("D22").Value = 67096
$ws.Range("E22").Value = 28156
$ws.Range("G22").Value = 44
$ws.Range("H22").Value = 8566

# Row 37: Singapur
$ws.Range("D37").Value = 37985
$ws.Range("E37").Value = 5650

# Row 43: Republica Dominicana (was Panama)
$ws.Range("A43").Value = 'Republica Dominicana'
$ws.Range("B43").Value = 31816
$ws.Range("C43").Value = 443
$ws.Range("D43").Value = 17280
$ws.Range("E43").Value = 13803
$ws.Range("G43").Value = 7
$ws.Range("H43").Value = 733

# Row 44: Panama (shifted down)
$ws.Range("A44").Value = 'Panama'
$ws.Range("B44").Value = 31686
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 15470
$ws.Range("E44").Value = 15612
$ws.Range("H44").Value = 604

# Row 45: Suiza (shifted down)
$ws.Range("A45").Value = 'Suiza'
$ws.Range("B45").Value = 31652
$ws.Range("C45").Value = 35
$ws.Range("D45").Value = 29100
$ws.Range("E45").Value = 590
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 1962

# Row 46: Bolivia (shifted down)
$ws.Range("A46").Value = 'Bolivia'
$ws.Range("B46").Value = 31524
$ws.Range("C46").Value = 848
$ws.Range("D46").Value = 8517
$ws.Range("E46").Value = 21993
$ws.Range("G46").Value = 44
$ws.Range("H46").Value = 1014

# Row 61: Moldavia
$ws.Range("B61").Value = 16357
$ws.Range("C61").Value = 107
$ws.Range("E61").Value = 6592
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 536

# Row 63: Argelia
$ws.Range("B63").Value = 13571
$ws.Range("C63").Value = 298
$ws.Range("D63").Value = 9671
$ws.Range("E63").Value = 2995
$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 905

# Row 69: Chequia
$ws.Range("B69").Value = 11724
$ws.Range("C69").Value = 121
$ws.Range("D69").Value = 7737
$ws.Range("E69").Value = 3639

# Row 79: Republica de Macedonia (was El Salvador)
$ws.Range("A79").Value = 'Republica de Macedonia'
$ws.Range("B79").Value = 6209
$ws.Range("C79").Value = 129
$ws.Range("D79").Value = 2427
$ws.Range("E79").Value = 3484
$ws.Range("H79").Value = 298

# Row 80: Kenia (shifted up into new slot)
$ws.Range("A80").Value = 'Kenia'
$ws.Range("B80").Value = 6190
$ws.Range("C80").Value = 120
$ws.Range("D80").Value = 2013
$ws.Range("E80").Value = 4033
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 144

# Row 81: El Salvador (shifted down)
$ws.Range("A81").Value = 'El Salvador'
$ws.Range("B81").Value = 6173
$ws.Range("C81").Value = 239
$ws.Range("D81").Value = 3648
$ws.Range("E81").Value = 2361
$ws.Range("G81").Value = 12
$ws.Range("H81").Value = 164

# Row 82: Tayikistan
$ws.Range("B82").Value = 5900
$ws.Range("C82").Value = 51
$ws.Range("D82").Value = 4506
$ws.Range("E82").Value = 1342

# Row 84: Etiopia
$ws.Range("B84").Value = 5846
$ws.Range("C84").Value = 157
$ws.Range("D84").Value = 2430
$ws.Range("E84").Value = 3313
$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 103

# Row 90: Republica de Yibuti
$ws.Range("B90").Value = 4656
$ws.Range("C90").Value = 13
$ws.Range("D90").Value = 4433
$ws.Range("E90").Value = 170
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 53

# Row 92: Luxemburgo
$ws.Range("B92").Value = 4256
$ws.Range("C92").Value = 14
$ws.Range("D92").Value = 3997
$ws.Range("E92").Value = 149

# Row 97: Grecia
$ws.Range("B97").Value = 3390
$ws.Range("C97").Value = 14
$ws.Range("E97").Value = 1825

# Row 104: Cuba
$ws.Range("B104").Value = 2340
$ws.Range("C104").Value = 8
$ws.Range("D104").Value = 2211
$ws.Range("E104").Value = 43

# Row 111: Sri Lanka
$ws.Range("B111").Value = 2039
$ws.Range("C111").Value = 2
$ws.Range("E111").Value = 350

# Row 129: Jordania
$ws.Range("B129").Value = 1128
$ws.Range("C129").Value = 7
$ws.Range("D129").Value = 867

# Row 144: Liberia
$ws.Range("B144").Value = 770
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 312
$ws.Range("E144").Value = 422
$ws.Range("G144").Value = 2
$ws.Range("H144").Value = 36

# Row 155: Montenegro
$ws.Range("B155").Value = 501
$ws.Range("C155").Value = 20
$ws.Range("E155").Value = 175

# Row 156: Surinam
$ws.Range("E156").Value = 280
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 13

